# Update "Linear Regression" results sheet with refreshed model output values
# (baseline re-run produced slightly different coefficients; values below come from the new run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = [double]"5.42754867889704E+35"
$ws.Range("C2").Value = [double]"1.507261706098355E+34"
$ws.Range("D2").Value = [double]"1.954424008956725E+37"
$ws.Range("E2").Value = [double]"77.67679566046903"
$ws.Range("F2").Value = [double]"1.828458224374995"
$ws.Range("G2").Value = [double]"42.4821276335261"
$ws.Range("H2").Value = [double]"0"

# Row 3
$ws.Range("B3").Value = [double]"68.66095035333535"
$ws.Range("C3").Value = [double]"64.27688385548889"
$ws.Range("D3").Value = [double]"73.16201468195662"
$ws.Range("E3").Value = [double]"0.5227203028711864"
$ws.Range("F3").Value = [double]"0.01343733084017988"
$ws.Range("G3").Value = [double]"38.90060526813589"
$ws.Range("H3").Value = [double]"2.238293274943599E-306"

# Row 4
$ws.Range("B4").Value = [double]"110.5638495180865"
$ws.Range("C4").Value = [double]"94.5825747049449"
$ws.Range("D4").Value = [double]"127.8576835110027"
$ws.Range("E4").Value = [double]"0.7446187442674873"
$ws.Range("F4").Value = [double]"0.04027158606163192"
$ws.Range("G4").Value = [double]"18.4899284355952"
$ws.Range("H4").Value = [double]"5.797234086427238E-75"

# Row 5
$ws.Range("B5").Value = [double]"4.067043413890148"
$ws.Range("C5").Value = [double]"-1.256454623025671"
$ws.Range("D5").Value = [double]"9.67754381881749"
$ws.Range("E5").Value = [double]"0.03986515366736176"
$ws.Range("F5").Value = [double]"0.02679045983757124"
$ws.Range("G5").Value = [double]"1.488035439072771"
$ws.Range("H5").Value = [double]"0.1367761636417744"

# Row 6
$ws.Range("B6").Value = [double]"387637.6638618933"
$ws.Range("C6").Value = [double]"189855.9516252367"
$ws.Range("D6").Value = [double]"791349.25279143"
$ws.Range("E6").Value = [double]"8.262914079875076"
$ws.Range("F6").Value = [double]"0.3640493750762665"
$ws.Range("G6").Value = [double]"22.69723462138628"
$ws.Range("H6").Value = [double]"5.488189412898144E-111"

# Row 7
$ws.Range("B7").Value = [double]"-32.97596724625498"
$ws.Range("C7").Value = [double]"-73.71301499039156"
$ws.Range("D7").Value = [double]"70.89144932114073"
$ws.Range("E7").Value = [double]"-0.4001189330967176"
$ws.Range("F7").Value = [double]"0.4775394401782307"
$ws.Range("G7").Value = [double]"-0.8378762033715631"
$ws.Range("H7").Value = [double]"0.4021222186931896"

# Row 8
$ws.Range("B8").Value = [double]"353.0136531720787"
$ws.Range("C8").Value = [double]"85.02603011760151"
$ws.Range("D8").Value = [double]"1009.14864157154"
$ws.Range("E8").Value = [double]"1.510752078493471"
$ws.Range("F8").Value = [double]"0.4568498703891006"
$ws.Range("G8").Value = [double]"3.306889585427173"
$ws.Range("H8").Value = [double]"0.0009470329022164"
